$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AffordableAMI")

# New source data: set the "Percent" (renter-occupied units) column (E) to 100
# for every data row (rows 2-62; row 1 is the header).
$ws.Range("E2:E62").Value = 100

# This sheet becomes the active tab in the saved workbook.
$ws.Activate()
